$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Use PO#" in H1
$ws.Range("H1").Value = "Use PO#"

# Update selection to H2 (single cell, no multi-range)
$ws.Range("H2").Select()
